$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.963.41"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.384.91"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.42"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.89"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.61"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.388"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.961.89"
$ws.Range("E12").Value = "  -0.46%  "
$ws.Range("E13").Value = "  +2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.81"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.433.33"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.069.83"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.08"
$ws.Range("E18").Value = "  -3.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.65"
$ws.Range("E19").Value = "  -3.75%  "
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.23"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.14"
$ws.Range("E22").Value = "  +2.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.521.47"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.15"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.97"
$ws.Range("E31").Value = "  -2.72%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("E33").Value = "  -5.76%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("E35").Value = "  -0.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.87"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.416.89"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.98"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0769"
$ws.Range("E40").Value = "  -2.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.76"
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.38"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.452.61"
$ws.Range("E47").Value = "  -3.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.94"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.72"
$ws.Range("E49").Value = "  -2.38%  "
$ws.Range("E50").Value = "  +9.29%  "
$ws.Range("E51").Value = "  +1.26%  "
